$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 135, pushing existing rows 135-218 down to 136-219.
$ws.Rows.Item(135).Insert()

# Populate the newly inserted row 135 with the new record's data.
$ws.Cells.Item(135, 1).Value = 8
$ws.Cells.Item(135, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(135, 3).Value = "Coquimbo"
$ws.Cells.Item(135, 4).Value = 45029
$ws.Cells.Item(135, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(135, 5).Value = 4
$ws.Cells.Item(135, 6).Value = 100112044
$ws.Cells.Item(135, 7).Value = "Perejil"
$ws.Cells.Item(135, 8).Value = "Sin especificar"
$ws.Cells.Item(135, 9).Value = "Primera"
$ws.Cells.Item(135, 10).Value = 2000
$ws.Cells.Item(135, 11).Value = 1800
$ws.Cells.Item(135, 12).Value = 2000
$ws.Cells.Item(135, 13).Value = 1900
$ws.Cells.Item(135, 14).Value = "$/atado 1 a 1,5 kilos"
$ws.Cells.Item(135, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(135, 16).Value = 1267
$ws.Cells.Item(135, 17).Value = 1.5
$ws.Cells.Item(135, 18).Value = "Hortaliza"
